# Commit: sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio
#
# This updates the "last_edited_time" (column D) timestamp for several Notion
# pages that were touched by this edit, and updates several numeric
# "properties" values on row 13 (the "Tháng 7" page) to reflect the corrected
# discount for phu sale and the new hourly-wage calculation strategy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- last_edited_time updates (column D) ---------------------------------
$newEditedTime = "2024-07-21T16:44:00.000Z"

$ws.Range("D4").Value  = $newEditedTime
$ws.Range("D5").Value  = $newEditedTime
$ws.Range("D6").Value  = $newEditedTime
$ws.Range("D8").Value  = $newEditedTime
$ws.Range("D12").Value = $newEditedTime
$ws.Range("D13").Value = $newEditedTime

# --- row 13 ("Tháng 7") numeric property updates --------------------------
# properties.Chi tiêu.number
$ws.Range("W13").Value = 64467000
# properties.Lũy kế.formula.number
$ws.Range("AA13").Value = 160813000
# properties.Tổng doanh thu.formula.number
$ws.Range("AE13").Value = 225280000
# properties.Đã thanh toán.number
$ws.Range("AH13").Value = 190580000
# properties.Số lượng đơn.number
$ws.Range("AK13").Value = 31
# properties.Đơn giá.number
$ws.Range("AQ13").Value = 216380000
